# Adds a new "M mambo marie" entry to the end of the used-icon list and
# relocates the hidden "_GoBack" bookmark so it again sits right after the
# newly-typed text (mirroring where Word leaves it after the last edit).

$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark (currently sitting right after
#    "% corpse", the paragraph that used to be last-edited).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Append a brand-new list paragraph after the current last item
#    ("r rifle ammo"). InsertParagraphAfter() clones the paragraph/list
#    formatting (ListParagraph style + numPr + rPr lang) of the paragraph
#    it's called on, which is exactly what every other bullet uses.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)

# 3. Type the full entry text in one go, then re-open a collapsed range
#    between "M mambo m" and "arie" for the bookmark. Splitting the text
#    into two runs this way (instead of collapsing at the paragraph's very
#    end) avoids the run/paragraph boundary and matches how the bookmark
#    ends up positioned between two runs in the final document.
$newPara.Range.Text = "M mambo marie"

$firstChunk = "M mambo m"
$splitPoint = $newPara.Range.Start + $firstChunk.Length
$bookmarkRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
